$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 471, shifting existing rows down
$ws.Range("A471:A473").EntireRow.Insert()

$ws.Range("A471").Value = 8
$ws.Range("B471").Value = "Terminal La Palmera de La Serena"
$ws.Range("C471").Value = "Coquimbo"
$ws.Range("D471").Value = 44995
$ws.Range("E471").Value = 4
$ws.Range("F471").Value = "Fruta"
$ws.Range("G471").Value = 100103
$ws.Range("H471").Value = "Frutos de hueso (carozo)"
$ws.Range("I471").Value = 100103006
$ws.Range("J471").Value = 'Nectarín'
$ws.Range("K471").Value = "Artic Pride"
$ws.Range("L471").Value = "Especial"
$ws.Range("M471").Value = 16
$ws.Range("N471").Value = 430000
$ws.Range("O471").Value = 440000
$ws.Range("P471").Value = 435000
$ws.Range("Q471").Value = '$/bins (420 kilos)'
$ws.Range("R471").Value = "Región de O'Higgins"
$ws.Range("S471").Value = 1036
$ws.Range("T471").Value = 420

$ws.Range("A472").Value = 8
$ws.Range("B472").Value = "Terminal La Palmera de La Serena"
$ws.Range("C472").Value = "Coquimbo"
$ws.Range("D472").Value = 44995
$ws.Range("E472").Value = 4
$ws.Range("F472").Value = "Fruta"
$ws.Range("G472").Value = 100103
$ws.Range("H472").Value = "Frutos de hueso (carozo)"
$ws.Range("I472").Value = 100103006
$ws.Range("J472").Value = 'Nectarín'
$ws.Range("K472").Value = "Artic Pride"
$ws.Range("L472").Value = "Primera"
$ws.Range("M472").Value = 16
$ws.Range("N472").Value = 400000
$ws.Range("O472").Value = 410000
$ws.Range("P472").Value = 405000
$ws.Range("Q472").Value = '$/bins (420 kilos)'
$ws.Range("R472").Value = "Región de O'Higgins"
$ws.Range("S472").Value = 964
$ws.Range("T472").Value = 420

$ws.Range("A473").Value = 8
$ws.Range("B473").Value = "Terminal La Palmera de La Serena"
$ws.Range("C473").Value = "Coquimbo"
$ws.Range("D473").Value = 44995
$ws.Range("E473").Value = 4
$ws.Range("F473").Value = "Fruta"
$ws.Range("G473").Value = 100103
$ws.Range("H473").Value = "Frutos de hueso (carozo)"
$ws.Range("I473").Value = 100103006
$ws.Range("J473").Value = 'Nectarín'
$ws.Range("K473").Value = "Artic Pride"
$ws.Range("L473").Value = "Segunda"
$ws.Range("M473").Value = 16
$ws.Range("N473").Value = 350000
$ws.Range("O473").Value = 360000
$ws.Range("P473").Value = 355000
$ws.Range("Q473").Value = '$/bins (420 kilos)'
$ws.Range("R473").Value = "Región de O'Higgins"
$ws.Range("S473").Value = 845
$ws.Range("T473").Value = 420
